$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: was =NOW() formatted as a date/time, now =A1+2 (a plain number).
# Re-point the cell's number format back to General before the row shift
# below, so the existing style slot is reused instead of a new one created.
$ws.Range("B1").Formula = "=A1+2"
$ws.Range("B1").NumberFormat = "General"

# Drop the old row 2 (A2="a", B2==1+2) entirely; old row 3 shifts up to
# become the new row 2.
$ws.Rows.Item(2).Delete()

# New row 2 (previously row 3): keep A2's =TODAY(), replace B2 (previously
# "=A3") with a text formula based on the date in A2.
$ws.Range("A2").Formula = "=TODAY()"
$ws.Range("B2").Formula = '=TEXT(A2,"YYYY-MM-DD")'

$ws.Range("B1").Select() | Out-Null
